$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New roster data (header + 17 players), replacing the old 18-player table.
$data = @(
    @("Oyuncu Adı", "Pozisyon", "Takım"),
    @("LaMelo Ball", "PG,SG", "Charlotte Hornets"),
    @("Collin Sexton", "PG,SG", "Utah Jazz"),
    @("Shaedon Sharpe", "SG,SF", "Portland Trail Blazers"),
    @("Cade Cunningham", "PG,SG", "Detroit Pistons"),
    @("Malik Monk", "PG,SG,SF", "Sacramento Kings"),
    @("Ausar Thompson", "SF,PF", "Detroit Pistons"),
    @("Onyeka Okongwu", "PF,C", "Atlanta Hawks"),
    @("Keon Ellis", "SG,SF", "Sacramento Kings"),
    @("Isaiah Hartenstein", "C", "Oklahoma City Thunder"),
    @("Aaron Nesmith", "SF,PF", "Indiana Pacers"),
    @("Coby White", "PG,SG", "Chicago Bulls"),
    @("Damian Lillard", "PG", "Milwaukee Bucks"),
    @("Derrick White", "PG,SG", "Boston Celtics"),
    @("Naz Reid", "PF,C", "Minnesota Timberwolves"),
    @("Gary Trent Jr.", "PG,SG,SF", "Milwaukee Bucks"),
    @("Anthony Davis", "PF,C", "Dallas Mavericks"),
    @("Andrew Wiggins", "SF,PF", "Miami Heat")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 1
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}

# The old table had one extra row (19) — remove it now that only 18 rows remain.
$ws.Range("A19:C19").Delete()
